$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.376.02"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "2.292.64"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.68"
$ws.Range("E5").Value = "  -1.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.80"
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  -2.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.49"
$ws.Range("E10").Value = "  -2.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.04"
$ws.Range("E11").Value = "  +3.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0782"
$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").Value = "2.648.62"
$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").Value = "2.281.76"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").Value = "42.338.82"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("E19").Value = "  -6.31%  "

$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -1.30%  "

$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.77"
$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("E23").Value = "  +6.33%  "

$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  -1.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.29"
$ws.Range("E27").Value = "  -4.18%  "

$ws.Range("E28").Value = "  +8.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.73"
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.04"
$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("E31").Value = "  -3.77%  "

$ws.Range("E33").Value = "  +0.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.55"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("E35").Value = "  -7.64%  "

$ws.Range("E36").Value = "  +0.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -2.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0999"
$ws.Range("E38").Value = "  -1.62%  "

$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("E40").Value = "  -1.25%  "

$ws.Range("E41").Value = "  -0.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.01"
$ws.Range("E42").Value = "  +9.93%  "

$ws.Range("D43").Value = "1.964.07"
$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.44"
$ws.Range("E44").Value = "  +4.07%  "

$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("E47").Value = "  -1.59%  "

$ws.Range("E48").Value = "  -1.50%  "

$ws.Range("D49").Value = "2.516.21"
$ws.Range("E49").Value = "  +0.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.11"
$ws.Range("E50").Value = "  -1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.01"
$ws.Range("E51").Value = "  -0.79%  "
